{"js": "// Turn the trailing empty paragraph (Courier New placeholder just before the\n// section break) into two new log-entry paragraphs:\n//   1) a bold \"List Paragraph\" numbered heading with the new timestamp, and\n//   2) the note paragraph (ind left=1080, Courier New) with the actual text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = paragraphs.items[paragraphs.items.length - 1];\n\nconst headingParaXml =\n  \"<w:p>\" +\n    \"<w:pPr>\" +\n      \"<w:pStyle w:val=\\\"ListParagraph\\\"/>\" +\n      \"<w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr>\" +\n      \"<w:rPr>\" +\n        \"<w:rFonts w:ascii=\\\"Times New Roman\\\" w:hAnsi=\\\"Times New Roman\\\" w:cs=\\\"Times New Roman\\\"/>\" +\n        \"<w:b/><w:sz w:val=\\\"24\\\"/><w:szCs w:val=\\\"24\\\"/><w:lang w:val=\\\"en-US\\\"/>\" +\n      \"</w:rPr>\" +\n    \"</w:pPr>\" +\n    \"<w:r>\" +\n      \"<w:rPr>\" +\n        \"<w:rFonts w:ascii=\\\"Times New Roman\\\" w:hAnsi=\\\"Times New Roman\\\" w:cs=\\\"Times New Roman\\\"/>\" +\n        \"<w:b/><w:sz w:val=\\\"24\\\"/><w:szCs w:val=\\\"24\\\"/><w:lang w:val=\\\"en-US\\\"/>\" +\n      \"</w:rPr>\" +\n      \"<w:t>0812506 \\u2013 2/06/2012 21:25</w:t>\" +\n    \"</w:r>\" +\n  \"</w:p>\";\n\nconst noteParaXml =\n  \"<w:p>\" +\n    \"<w:pPr>\" +\n      \"<w:ind w:left=\\\"1080\\\"/>\" +\n      \"<w:rPr>\" +\n        \"<w:rFonts w:ascii=\\\"Courier New\\\" w:hAnsi=\\\"Courier New\\\" w:cs=\\\"Courier New\\\"/>\" +\n        \"<w:noProof/><w:sz w:val=\\\"20\\\"/><w:szCs w:val=\\\"20\\\"/><w:lang w:val=\\\"en-US\\\"/>\" +\n      \"</w:rPr>\" +\n    \"</w:pPr>\" +\n    \"<w:r>\" +\n      \"<w:rPr>\" +\n        \"<w:rFonts w:ascii=\\\"Courier New\\\" w:hAnsi=\\\"Courier New\\\" w:cs=\\\"Courier New\\\"/>\" +\n        \"<w:noProof/><w:sz w:val=\\\"20\\\"/><w:szCs w:val=\\\"20\\\"/><w:lang w:val=\\\"en-US\\\"/>\" +\n      \"</w:rPr>\" +\n      \"<w:t>Ch\\u1EC9nh l\\u1EA1i set @MonHoc trong usp_CapNhatSoLuongSVNhom_Error v\\u00E0 usp_CapNhatSoLuongSVNhom_Fix</w:t>\" +\n    \"</w:r>\" +\n  \"</w:p>\";\n\n// Office.js's insertOoxml requires a full \"Flat OPC\" package envelope (not a\n// bare <w:p> fragment like the COM Range.InsertXML accepts), so wrap the two\n// target paragraphs in a minimal one and swap them in for the placeholder\n// paragraph in one shot.\nconst flatOpc =\n  \"<?xml version=\\\"1.0\\\" standalone=\\\"yes\\\"?>\" +\n  \"<pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\">\" +\n    \"<pkg:part pkg:name=\\\"/_rels/.rels\\\" pkg:contentType=\\\"application/vnd.openxmlformats-package.relationships+xml\\\" pkg:padding=\\\"512\\\">\" +\n      \"<pkg:xmlData>\" +\n        \"<Relationships xmlns=\\\"http://schemas.openxmlformats.org/package/2006/relationships\\\">\" +\n          \"<Relationship Id=\\\"rId1\\\" Type=\\\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\\\" Target=\\\"word/document.xml\\\"/>\" +\n        \"</Relationships>\" +\n      \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"<pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\">\" +\n      \"<pkg:xmlData>\" +\n        \"<w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\">\" +\n          \"<w:body>\" + headingParaXml + noteParaXml + \"</w:body>\" +\n        \"</w:document>\" +\n      \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ntarget.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Turn the trailing empty paragraph (Courier New placeholder just before the\n# section break) into two new log-entry paragraphs:\n#   1) a bold \"List Paragraph\" numbered heading with the new timestamp, and\n#   2) the note paragraph (ind left=1080, Courier New) with the actual text.\n$d = $word.ActiveDocument\n\n$target = $d.Paragraphs.Last\n$r = $target.Range\n$r.Collapse(1)\n\n$wNs = \"xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'\"\n\n$headingXml = \"<w:p $wNs>\" +\n  \"<w:pPr>\" +\n    \"<w:pStyle w:val='ListParagraph'/>\" +\n    \"<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>\" +\n    \"<w:rPr>\" +\n      \"<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>\" +\n      \"<w:b/><w:sz w:val='24'/><w:szCs w:val='24'/><w:lang w:val='en-US'/>\" +\n    \"</w:rPr>\" +\n  \"</w:pPr>\" +\n  \"<w:r>\" +\n    \"<w:rPr>\" +\n      \"<w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman' w:cs='Times New Roman'/>\" +\n      \"<w:b/><w:sz w:val='24'/><w:szCs w:val='24'/><w:lang w:val='en-US'/>\" +\n    \"</w:rPr>\" +\n    \"<w:t>0812506 \u2013 2/06/2012 21:25</w:t>\" +\n  \"</w:r>\" +\n\"</w:p>\"\n\n$noteXml = \"<w:p $wNs>\" +\n  \"<w:pPr>\" +\n    \"<w:ind w:left='1080'/>\" +\n    \"<w:rPr>\" +\n      \"<w:rFonts w:ascii='Courier New' w:hAnsi='Courier New' w:cs='Courier New'/>\" +\n      \"<w:noProof/><w:sz w:val='20'/><w:szCs w:val='20'/><w:lang w:val='en-US'/>\" +\n    \"</w:rPr>\" +\n  \"</w:pPr>\" +\n  \"<w:r>\" +\n    \"<w:rPr>\" +\n      \"<w:rFonts w:ascii='Courier New' w:hAnsi='Courier New' w:cs='Courier New'/>\" +\n      \"<w:noProof/><w:sz w:val='20'/><w:szCs w:val='20'/><w:lang w:val='en-US'/>\" +\n    \"</w:rPr>\" +\n    \"<w:t>Ch\u1ec9nh l\u1ea1i set @MonHoc trong usp_CapNhatSoLuongSVNhom_Error v\u00e0 usp_CapNhatSoLuongSVNhom_Fix</w:t>\" +\n  \"</w:r>\" +\n\"</w:p>\"\n\n$r.InsertXML($headingXml + $noteXml)\n"}
